# Auto-generated edit script: update Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet to reflect refreshed market data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '26.660.20'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  +1.02%  '

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '1.854.14'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  +0.63%  '

$ws.Range('E4').Value = '  +0.07%  '

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '264.92'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +1.64%  '

$ws.Range('E6').Value = '  +0.02%  '

$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.5278'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  +0.62%  '

$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.3249'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +0.47%  '

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.06803'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +0.89%  '

$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '18.99'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  +0.90%  '

$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.7836'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +1.59%  '

$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.07801'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  +1.66%  '

$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '1.860.61'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +0.91%  '

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '88.64'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -0.30%  '

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '5.035'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  +0.18%  '

$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '1.002'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +0.13%  '

$ws.Range('E17').Value = '  -0.73%  '

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '0.000007983'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +1.25%  '

$ws.Range('E19').Value = '  +0.04%  '

$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '26.684.37'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +0.92%  '

$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '4.653'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +2.73%  '

$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '9.500'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +0.65%  '

$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '6.018'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +1.59%  '

$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '143.21'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -0.71%  '

$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '2.179'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -6.85%  '

$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '1.695'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +2.51%  '

$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '17.03'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +0.73%  '

$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '112.18'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  +0.61%  '

$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '4.199'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +0.19%  '

$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '4.117'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +0.02%  '

$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '0.08732'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -0.73%  '

$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '0.04852'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +0.15%  '

$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '0.7226'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +5.23%  '

$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '1.133'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +0.08%  '

$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '2.879'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  +0.72%  '

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '3.116'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +0.19%  '

$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '2.273'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +2.70%  '

$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.01798'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +0.49%  '

$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '0.4883'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -1.04%  '

$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.9057'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +0.47%  '

$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '111.25'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -1.39%  '

$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '5.979'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -3.18%  '

$ws.Range('E43').Value = '  +0.06%  '

$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '7.702'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -0.40%  '

$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.4216'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +0.57%  '

$ws.Range('E46').Value = '  +0.21%  '

$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '9.033'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -0.80%  '

$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '0.1241'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -1.66%  '

$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '35.15'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -0.68%  '

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.8897'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +3.36%  '

$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '60.20'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +1.65%  '

